# TestCases.xlsx update:
#  - Row 8 (Test Case 7, "Test invalid input while playing the maze" / console-capture case):
#    mark the Actual Outcome / Fail-Pass / Justification columns as "Will not be doable as
#    console output cannot be captured while the program is asking for an input."
#  - Remove the screenshot picture that used to sit in column F (no longer applicable now
#    that the row records "not doable" instead of an actual-outcome screenshot), and
#    shrink column F back down from the width that was only needed to host that picture.
#  - Leave the selection on F9, matching where the user ended up after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Test Case 7 (row 8): Actual Outcome / Fail-Pass / Justification.
$ws.Range("F8").Value = "-"
$ws.Range("G8").Value = "-"
$ws.Range("H8").Value = "Will not be doable as console output cannot be captured while the program is asking for an input."

# Delete the embedded screenshot picture anchored over column F / row 8.
for ($i = $ws.Shapes.Count; $i -ge 1; $i--) {
    $ws.Shapes.Item($i).Delete()
}

# Column F no longer needs to be wide enough to host the picture.
$ws.Columns.Item(6).ColumnWidth = 33.67

# Restore the active selection to F9.
$ws.Range("F9").Select()
